# Adds the "WYSIWYG Canvas Deep-Nesting Update" status section to the end
# of the document. The new section mirrors the structure of the existing
# "Visual Drag-Drop Canvas Foundation Update" block that precedes it:
# a blank line, a "---" separator, a title line, an "Updated:" line, a
# blank line, a tab-delimited header row, and two tab-delimited data rows.

$d = $word.ActiveDocument

$wordNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# Every run in this document shares the same character formatting
# (Helvetica Light, 12pt / half-point size 24), matching the rest of the
# report's body text.
$runProps = '<w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr>'

function Escape-Xml($text) {
    $escaped = $text.Replace("&", "&amp;")
    $escaped = $escaped.Replace("<", "&lt;")
    $escaped = $escaped.Replace(">", "&gt;")
    return $escaped
}

# NOTE: this interpreter shares one flat variable scope between callers and
# callees, so every loop at every nesting level needs its own, uniquely
# named counter ($ci here) - reusing a name like $i that an outer loop is
# also using will stomp on the outer loop's counter.
function New-ReportParagraphXml($cells, $isFirst) {
    $run = $runProps
    for ($ci = 0; $ci -lt $cells.Count; $ci++) {
        if ($ci -gt 0) {
            $run += '<w:tab/>'
        }
        $escaped = Escape-Xml $cells[$ci]
        $run += '<w:t xml:space="preserve">' + $escaped + '</w:t>'
    }

    if ($isFirst) {
        return '<w:p xmlns:w="' + $wordNs + '"><w:pPr/><w:r>' + $run + '</w:r></w:p>'
    }
    return '<w:p><w:pPr/><w:r>' + $run + '</w:r></w:p>'
}

# Each entry is one paragraph; each element of the inner array is one
# tab-separated cell within that paragraph (a single-element array is a
# plain text line with no tabs).
$rows = @(
    , @("")
    , @("---")
    , @("WYSIWYG Canvas Deep-Nesting Update")
    , @("Updated: 2026-02-18")
    , @("")
    , @("Module Name", "Developed", "Partial Developed", "Need To Develop")
    , @(
        "Visual Drag-Drop Canvas",
        "Added freeform nested canvas rendering with recursive drop targets, node-level drag/reparent as child, advanced inspector fields (type/title/style settings), and JSON-backed node identity model",
        "Reparent currently supports drop-as-child (no before/after drop zones yet), no multi-select/group operations",
        "Pixel-perfect freeform absolute positioning mode, container grid snapping, full drag handles by breakpoint"
      )
    , @(
        "Editor UX",
        "Undo/redo timeline maintained with deep snapshots; selected-node inspector edits settings JSON + style fields in real time",
        "No timeline scrubber UI and no collaborative conflict resolution",
        "Operational transform/CRDT collaboration, timeline diff viewer, visual history playback"
      )
)

$sectionXml = ""
for ($pi = 0; $pi -lt $rows.Count; $pi++) {
    $isFirst = $false
    if ($pi -eq 0) {
        $isFirst = $true
    }
    $rowCells = $rows[$pi]
    $sectionXml += New-ReportParagraphXml $rowCells $isFirst
}

# Append after the very last paragraph, ahead of the section mark.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertXML($sectionXml)

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
